$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

for ($r = $startRow; $r -lt ($startRow + $rowCount); $r++) {
    for ($c = $startCol; $c -lt ($startCol + $colCount); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -is [double] -or $val -is [int]) {
            $dval = [double]$val
            if ($dval -ge 0) {
                $rounded = [math]::Floor($dval + 0.5)
            } else {
                $rounded = [math]::Ceiling($dval - 0.5)
            }
            $cell.Value = $rounded
        }
    }
}
